# Add a new "25-ago" day column (BC) to the right of the existing "24-ago"
# column (BB) on Sheet1, mirroring the formatting of column BB and filling
# in the new daily values for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number format, alignment, style) of column BB into
# the new column BC so the new cells match the look of the rest of the
# table (header style for row 1, centered-integer style for the data rows).
$ws.Range("BB1:BB11").Copy() | Out-Null
$ws.Range("BC1:BC11").PasteSpecial(-4122) | Out-Null

# Header for the new date column.
$ws.Range("BC1").Value = "25-ago"

# New day's values for each product row.
$ws.Range("BC2").Value = 16
$ws.Range("BC3").Value = 9
$ws.Range("BC4").Value = 9
$ws.Range("BC5").Value = 11
$ws.Range("BC6").Value = 12
$ws.Range("BC7").Value = 14
$ws.Range("BC8").Value = 11
$ws.Range("BC9").Value = 16
$ws.Range("BC10").Value = 25
$ws.Range("BC11").Value = 14

# Leave the selection where the author ended up after the edit.
$ws.Range("BI8").Select() | Out-Null
